$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Range("A2").Value = "Ayati Arvind"
$ws.Range("C3").Select()
